$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows (2-36) down one row, to 3-37, so a new
# "Total" summary row can be inserted at the top of the table, under the
# header row. Shifting the values manually (bottom-up, so a source row is
# always read before it gets overwritten) keeps every row's original
# formatting in place instead of minting new duplicated styles the way a
# native row-Insert would.
for ($r = 36; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Range("A$dest").Value = $ws.Range("A$r").Value2
    $ws.Range("B$dest").Value = $ws.Range("B$r").Value2
}

# Row 37 is brand-new, so its column-A cell has no style yet (column B
# already gets one from the column-level default). Give it the same
# look as the rest of column A by copying the format from row 36.
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "Total" row. Write column B first so the shared-string
# table registers the single-space placeholder before "Total".
$ws.Range("B2").Value = " "
$ws.Range("A2").Value = "Total"

# Reset the active selection to A3, matching the saved view state.
$ws.Range("A3").Select()
